# [IMP] re-arrange data mapping to templates
#
# Rebuild the worksheet header/label block to match the new template
# layout: a title row, three label rows (center / warehouse / report
# date), a spacer row, and a re-ordered column header row with a new
# "currency" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start from a clean slate: delete the old rows (1-9) entirely so no
# stray values, row heights or per-cell styles survive.
$ws.Rows("1:9").Delete()

function Set-LeftCell($addr) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Tahoma"
    $r.Font.Size = 10
    $r.WrapText = $true
    $r.HorizontalAlignment = -4131
}

function Set-RightCell($addr) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Tahoma"
    $r.Font.Size = 10
    $r.WrapText = $true
    $r.HorizontalAlignment = -4152
}

# ---------------------------------------------------------------
# Row 1 - report title (merged A1:B1), bold Tahoma 10
# ---------------------------------------------------------------
$ws.Range("A1:B1").Merge()
$ws.Range("A1").Value = "รายงานตรวจสอบจำนวน Stock คงเหลือ"

$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").Font.Name = "Tahoma"
$ws.Range("A1:B1").Font.Size = 10
$ws.Range("A1:B1").WrapText = $true
$ws.Range("A1:B1").HorizontalAlignment = -4131

Set-LeftCell "C1"
Set-RightCell "D1"
Set-LeftCell "E1"
Set-RightCell "F1"

$ws.Rows(1).RowHeight = 14

# ---------------------------------------------------------------
# Row 2 - center label, bold Tahoma 10
# ---------------------------------------------------------------
$ws.Range("A2").Value = "ศูนย์"

$ws.Range("A2:B2").Font.Bold = $true
$ws.Range("A2:B2").Font.Name = "Tahoma"
$ws.Range("A2:B2").Font.Size = 10
$ws.Range("A2:B2").WrapText = $true
$ws.Range("A2:B2").HorizontalAlignment = -4131

Set-LeftCell "C2"
Set-RightCell "D2"
Set-LeftCell "E2"
Set-RightCell "F2"

$ws.Rows(2).RowHeight = 14

# ---------------------------------------------------------------
# Row 3 - warehouse label, bold Tahoma 10
# ---------------------------------------------------------------
$ws.Range("A3").Value = "คลังสินค้า"

$ws.Range("A3").Font.Bold = $true
$ws.Range("A3").Font.Name = "Tahoma"
$ws.Range("A3").Font.Size = 10
$ws.Range("A3").WrapText = $true
$ws.Range("A3").HorizontalAlignment = -4131

Set-LeftCell "B3"
Set-LeftCell "C3"
Set-RightCell "D3"
Set-LeftCell "E3"
Set-RightCell "F3"

$ws.Rows(3).RowHeight = 14

# ---------------------------------------------------------------
# Row 4 - report-run-date label, bold Tahoma 10
# ---------------------------------------------------------------
$ws.Range("A4").Value = "วันที่เรียกรายงาน"

$ws.Range("A4").Font.Bold = $true
$ws.Range("A4").Font.Name = "Tahoma"
$ws.Range("A4").Font.Size = 10
$ws.Range("A4").WrapText = $true
$ws.Range("A4").HorizontalAlignment = -4131

Set-LeftCell "B4"
Set-LeftCell "C4"
Set-RightCell "D4"
Set-LeftCell "E4"
Set-RightCell "F4"

$ws.Rows(4).RowHeight = 14

# ---------------------------------------------------------------
# Row 5 - spacer row (bold Tahoma 9 / Tahoma 9)
# ---------------------------------------------------------------
$ws.Range("A5").Font.Bold = $true
$ws.Range("A5").Font.Name = "Tahoma"
$ws.Range("A5").Font.Size = 9
$ws.Range("A5").WrapText = $true
$ws.Range("A5").HorizontalAlignment = -4131

$ws.Range("B5").Font.Name = "Tahoma"
$ws.Range("B5").Font.Size = 9
$ws.Range("B5").WrapText = $true
$ws.Range("B5").HorizontalAlignment = -4131

$ws.Rows(5).RowHeight = 14

# ---------------------------------------------------------------
# Row 6 - column headers (re-ordered + new "currency" column)
# ---------------------------------------------------------------
$ws.Range("A6").Value = "ชื่อสินค้า"
$ws.Range("B6").Value = "ศูนย์"
$ws.Range("C6").Value = "คลังวัสดุ"
$ws.Range("D6").Value = "ยอดคงเหลือ"
$ws.Range("E6").Value = "หน่วยนับ"
$ws.Range("F6").Value = "มูลค่าสินค้าคงคลัง"
$ws.Range("G6").Value = "สกุลเงิน"

$ws.Range("A6:G6").Font.Bold = $true
$ws.Range("A6:G6").Font.Name = "Tahoma"
$ws.Range("A6:G6").Font.Size = 9
$ws.Range("A6:G6").WrapText = $true
$ws.Range("A6:G6").HorizontalAlignment = -4108
$ws.Range("A6:G6").Borders.LineStyle = 1

$ws.Rows(6).RowHeight = 16.5

# ---------------------------------------------------------------
# Final touches
# ---------------------------------------------------------------
$ws.Range("B9").Select()
